# Updates cryptos list figures (price/volume columns) to the latest scraped
# values, matching the GitHub Actions "cryptos list" refresh commit.
# Note: some "Price" column values are plain numeric-looking strings
# (e.g. "188.02"); prefixing them with a leading apostrophe forces Excel to
# keep them stored as text, matching the original inline-string cell type
# instead of letting Excel auto-convert them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.908.49'
$ws.Range('E2').Value = '  +4.07%  '
$ws.Range('D3').Value = '3.531.36'
$ws.Range('E3').Value = '  +4.40%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '''188.02'
$ws.Range('E5').Value = '  +7.35%  '
$ws.Range('D6').Value = '''560.84'
$ws.Range('E6').Value = '  +6.55%  '
$ws.Range('D7').Value = '''0.629'
$ws.Range('E7').Value = '  +5.48%  '
$ws.Range('D8').Value = '3.530.01'
$ws.Range('E8').Value = '  +4.37%  '
$ws.Range('E9').Value = '  -0.04%  '
$ws.Range('D10').Value = '''0.635'
$ws.Range('E10').Value = '  +3.94%  '
$ws.Range('E11').Value = '  +15.99%  '
$ws.Range('D12').Value = '''54.92'
$ws.Range('E12').Value = '  +2.50%  '
$ws.Range('E13').Value = '  +6.30%  '
$ws.Range('D14').Value = '''9.35'
$ws.Range('E14').Value = '  +2.48%  '
$ws.Range('D15').Value = '4.090.48'
$ws.Range('E15').Value = '  +4.70%  '
$ws.Range('D16').Value = '3.533.59'
$ws.Range('E16').Value = '  +4.99%  '
$ws.Range('D17').Value = '''18.68'
$ws.Range('E17').Value = '  +5.85%  '
$ws.Range('E18').Value = '  +2.81%  '
$ws.Range('D19').Value = '66.919.22'
$ws.Range('E19').Value = '  +4.31%  '
$ws.Range('D20').Value = '''12.11'
$ws.Range('E20').Value = '  +7.05%  '
$ws.Range('E21').Value = '  +3.22%  '
$ws.Range('D22').Value = '''422.64'
$ws.Range('E22').Value = '  +12.60%  '
$ws.Range('D23').Value = '''4.09'
$ws.Range('E23').Value = '  +10.10%  '
$ws.Range('D24').Value = '''86.17'
$ws.Range('E24').Value = '  +5.70%  '
$ws.Range('E25').Value = '  +1.04%  '
$ws.Range('D26').Value = '''11.06'
$ws.Range('E26').Value = '  -4.94%  '
$ws.Range('D27').Value = '''2.91'
$ws.Range('E27').Value = '  +7.16%  '
$ws.Range('B28').Value = 'LEO'
$ws.Range('C28').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D28').Value = '''6.11'
$ws.Range('E28').Value = '  -0.93%  '
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').Value = '''12.28'
$ws.Range('E29').Value = '  +8.20%  '
$ws.Range('D30').Value = '''9.09'
$ws.Range('E30').Value = '  +10.02%  '
$ws.Range('D31').Value = '''30.38'
$ws.Range('E31').Value = '  +4.67%  '
$ws.Range('D32').Value = '''636.84'
$ws.Range('E32').Value = '  +0.81%  '
$ws.Range('E33').Value = '  +2.78%  '
$ws.Range('E34').Value = '  +4.35%  '
$ws.Range('D35').Value = '''0.111'
$ws.Range('E35').Value = '  +4.87%  '
$ws.Range('D36').Value = '''60.36'
$ws.Range('E36').Value = '  +3.96%  '
$ws.Range('D37').Value = '0.0₃0831'
$ws.Range('E37').Value = '  +11.65%  '
$ws.Range('B38').Value = 'InjectiveProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D38').Value = '''38.34'
$ws.Range('E38').Value = '  +5.21%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').Value = '''0.148'
$ws.Range('E39').Value = '  +17.88%  '
$ws.Range('E40').Value = '  -0.18%  '
$ws.Range('E41').Value = '  +1.55%  '
$ws.Range('E42').Value = '  +12.73%  '
$ws.Range('D43').Value = '3.124.19'
$ws.Range('E43').Value = '  +4.90%  '
$ws.Range('E44').Value = '  +0.22%  '
$ws.Range('E45').Value = '  -2.34%  '
$ws.Range('D46').Value = '''2.86'
$ws.Range('E46').Value = '  +9.33%  '
$ws.Range('D47').Value = '''3.35'
$ws.Range('E47').Value = '  +10.60%  '
$ws.Range('D48').Value = '''0.0418'
$ws.Range('E48').Value = '  +5.15%  '
$ws.Range('E49').Value = '  +2.21%  '
$ws.Range('E50').Value = '  +5.79%  '
$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D51').Value = '''8.53'
$ws.Range('E51').Value = '  +8.02%  '
